$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the metadata "type" tags in row 3 (sdmx/iaest dimension-measure markers) ---
$ws.Range("B3").Value = "iaest-measure:orden"
$ws.Range("D3").Value = "iaest-measure:siglas"

# --- Update row 4 (dim / medida classification) ---
$ws.Range("B4").Value = "medida"
$ws.Range("D4").Value = "medida"

# --- Update row 5 (xsd type / URI markers) ---
$ws.Range("B5").Value = "xsd:int"
$ws.Range("D5").Value = "xsd:string"

# --- Normalize formatting: columns B (rows 3-5) and D (rows 1-5) were carrying a
#     leftover explicit "Arial" font override; bring them back in line with the
#     sheet's normal (unstyled) cell format, same as column A. ---
$ws.Range("A1").Copy()
$ws.Range("D1:D5").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("B3:B5").PasteSpecial(-4122)

# --- Remove row 6 entirely (it only held the now-removed "mapping-siglas.xlsx" note) ---
$ws.Range("D6").EntireRow.Delete()
